$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsAbout = $wb.Worksheets.Item("About")
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")

# --- About sheet: update the "last updated" date in C1 ---
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity sheet: update capacity credit multipliers for the two
#     hydrogen rows (hydrogen combustion turbine / hydrogen combined cycle)
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Make RAF-capacity the active / selected sheet, with B25 selected and
#     zoomed to 80%, matching the saved workbook view state ---
$wsCapacity.Activate()
$wsCapacity.Columns.Item(1).ColumnWidth = 29
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
